# Applies the cryptos.xlsx price/volume update described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that writes a value as plain text, preventing Excel's automatic
# number/date conversion for strings such as "42.45" or "231.03" while
# leaving the cell's style/format exactly as it was before the write.
function Set-TextValue([string]$addr, [string]$val) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue 'D2' '35.683.99'
Set-TextValue 'E2' '  +3.45%  '
Set-TextValue 'D3' '1.863.54'
Set-TextValue 'E3' '  +2.91%  '
Set-TextValue 'E4' '  +0.41%  '
Set-TextValue 'D5' '231.03'
Set-TextValue 'E5' '  +2.37%  '
Set-TextValue 'D6' '0.611'
Set-TextValue 'E6' '  +2.94%  '
Set-TextValue 'D8' '42.45'
Set-TextValue 'E8' '  +10.64%  '
Set-TextValue 'E9' '  +7.11%  '
Set-TextValue 'E10' '  +3.11%  '
Set-TextValue 'E11' '  +4.15%  '
Set-TextValue 'D12' '2.132.88'
Set-TextValue 'E12' '  +2.89%  '
Set-TextValue 'D13' '11.63'
Set-TextValue 'E13' '  +3.65%  '
Set-TextValue 'D14' '1.858.43'
Set-TextValue 'E14' '  +2.61%  '
Set-TextValue 'E15' '  +7.28%  '
Set-TextValue 'E16' '  +6.69%  '
Set-TextValue 'D17' '35.694.43'
Set-TextValue 'E17' '  +3.56%  '
Set-TextValue 'D18' '70.43'
Set-TextValue 'E18' '  +3.04%  '
Set-TextValue 'D19' '248.96'
Set-TextValue 'E19' '  +2.24%  '
Set-TextValue 'E20' '  +4.13%  '
Set-TextValue 'D21' '12.28'
Set-TextValue 'E21' '  +9.45%  '
Set-TextValue 'D22' '4.76'
Set-TextValue 'E22' '  +15.23%  '
Set-TextValue 'E23' '  +0.32%  '
Set-TextValue 'E24' '  +0.51%  '
Set-TextValue 'D25' '170.61'
Set-TextValue 'E25' '  -0.08%  '
Set-TextValue 'E26' '  +2.89%  '
Set-TextValue 'E27' '  +1.41%  '
Set-TextValue 'E28' '  +1.45%  '
Set-TextValue 'E29' '  +16.62%  '
Set-TextValue 'E30' '  +0.48%  '
Set-TextValue 'D31' '3.349.30'
Set-TextValue 'E31' '  +37.85%  '
Set-TextValue 'B32' 'Hedera'
Set-TextValue 'C32' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D32' '0.0547'
Set-TextValue 'E32' '  +5.77%  '
Set-TextValue 'B33' 'InternetComputer(DFINITY)'
Set-TextValue 'C33' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D33' '4.09'
Set-TextValue 'E33' '  +6.01%  '
Set-TextValue 'D34' '3.95'
Set-TextValue 'E34' '  +4.06%  '
Set-TextValue 'E35' '  +3.91%  '
Set-TextValue 'D36' '101.87'
Set-TextValue 'E36' '  +24.21%  '
Set-TextValue 'D37' '0.695'
Set-TextValue 'E37' '  +7.89%  '
Set-TextValue 'D38' '1.372.51'
Set-TextValue 'E38' '  +1.18%  '
Set-TextValue 'D39' '2.50'
Set-TextValue 'E39' '  +6.68%  '
Set-TextValue 'E40' '  +3.19%  '
Set-TextValue 'D41' '0.0196'
Set-TextValue 'E41' '  +4.65%  '
Set-TextValue 'E42' '  +5.92%  '
Set-TextValue 'E43' '  +3.98%  '
Set-TextValue 'D44' '14.86'
Set-TextValue 'E44' '  +8.21%  '
Set-TextValue 'D45' '2.48'
Set-TextValue 'E45' '  +1.33%  '
Set-TextValue 'E46' '  +1.43%  '
Set-TextValue 'E47' '  +8.60%  '
Set-TextValue 'D48' '0.0523'
Set-TextValue 'E48' '  +2.61%  '
Set-TextValue 'D49' '2.031.35'
Set-TextValue 'E49' '  +2.91%  '
Set-TextValue 'D50' '104.43'
Set-TextValue 'E50' '  +1.79%  '
Set-TextValue 'E51' '  +0.42%  '
